# "Generate Report for Archive"
#
# 1. The status text "Ready for handoff" (used on the Overview sheet's
#    zh-cn/de-de columns and on the per-language Status column of the
#    zh-cn / de-de sheets) moves to "In Translation".
# 2. The two columns that held that status text are narrowed from
#    ~17.22 chars to ~13.41 chars:
#       - Overview sheet: columns E (zh-cn) and F (de-de)
#       - zh-cn sheet:    column C (Status)
#       - de-de sheet:    column C (Status)

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Update every occurrence of the old status text in place so all cells
# that shared the string now share the new one.
$overview.Cells.Replace($oldStatus, $newStatus)
$zhcn.Cells.Replace($oldStatus, $newStatus)
$dede.Cells.Replace($oldStatus, $newStatus)

# Narrow the status columns to match the new (shorter) content width.
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
